$d = $word.ActiveDocument

# The three bulleted list items ("Properly structured...", "Awesome",
# "...Very Awesome") were demoted from list level 1 (w:ilvl=0) to
# list level 2 (w:ilvl=1), i.e. indented one level deeper.
$targets = @(
    "Properly structured and follow all good OOP practices",
    "Awesome",
    "...Very Awesome"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    if ($targets -contains $text) {
        $p.Range.ListFormat.ListLevelNumber = 2
    }
}
